# Auto-generated edit script applying the cryptos.xlsx data refresh diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText {
    param($Cell, [string]$Text)
    # Force text interpretation so numeric-looking strings (e.g. "228.21")
    # are not silently converted into floating point numbers, then strip the
    # temporary formatting so the cell keeps its original (default) style.
    $Cell.NumberFormat = "@"
    $Cell.Value = $Text
    $Cell.ClearFormats()
}

Set-CellText $ws.Cells.Item(2, 4) "37.873.27"
Set-CellText $ws.Cells.Item(2, 5) "  -0.65%  "
Set-CellText $ws.Cells.Item(3, 4) "2.035.47"
Set-CellText $ws.Cells.Item(3, 5) "  -1.05%  "
Set-CellText $ws.Cells.Item(4, 5) "  +0.02%  "
Set-CellText $ws.Cells.Item(5, 4) "228.21"
Set-CellText $ws.Cells.Item(5, 5) "  -0.79%  "
Set-CellText $ws.Cells.Item(6, 4) "0.609"
Set-CellText $ws.Cells.Item(6, 5) "  -1.29%  "
Set-CellText $ws.Cells.Item(7, 4) "60.38"
Set-CellText $ws.Cells.Item(7, 5) "  +1.79%  "
Set-CellText $ws.Cells.Item(8, 5) "  +0.03%  "
Set-CellText $ws.Cells.Item(9, 4) "0.378"
Set-CellText $ws.Cells.Item(9, 5) "  -1.87%  "
Set-CellText $ws.Cells.Item(10, 5) "  +1.31%  "
Set-CellText $ws.Cells.Item(11, 5) "  -0.45%  "
Set-CellText $ws.Cells.Item(12, 4) "2.337.15"
Set-CellText $ws.Cells.Item(12, 5) "  -1.13%  "
Set-CellText $ws.Cells.Item(13, 4) "14.51"
Set-CellText $ws.Cells.Item(13, 5) "  -1.63%  "
Set-CellText $ws.Cells.Item(14, 4) "21.08"
Set-CellText $ws.Cells.Item(14, 5) "  -0.51%  "
Set-CellText $ws.Cells.Item(15, 5) "  +0.56%  "
Set-CellText $ws.Cells.Item(16, 5) "  -1.97%  "
Set-CellText $ws.Cells.Item(17, 4) "2.030.84"
Set-CellText $ws.Cells.Item(17, 5) "  -1.57%  "
Set-CellText $ws.Cells.Item(18, 4) "37.800.11"
Set-CellText $ws.Cells.Item(18, 5) "  -0.59%  "
Set-CellText $ws.Cells.Item(19, 5) "  -0.08%  "
Set-CellText $ws.Cells.Item(20, 4) "5.93"
Set-CellText $ws.Cells.Item(20, 5) "  -5.27%  "
Set-CellText $ws.Cells.Item(21, 5) "  -1.29%  "
Set-CellText $ws.Cells.Item(22, 4) "223.83"
Set-CellText $ws.Cells.Item(22, 5) "  -0.56%  "
Set-CellText $ws.Cells.Item(23, 5) "  +0.00%  "
Set-CellText $ws.Cells.Item(24, 4) "2.43"
Set-CellText $ws.Cells.Item(24, 5) "  -0.12%  "
Set-CellText $ws.Cells.Item(25, 4) "2.27"
Set-CellText $ws.Cells.Item(25, 5) "  +0.92%  "
Set-CellText $ws.Cells.Item(26, 2) "Monero"
Set-CellText $ws.Cells.Item(26, 3) "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-CellText $ws.Cells.Item(26, 4) "167.46"
Set-CellText $ws.Cells.Item(26, 5) "  +0.66%  "
Set-CellText $ws.Cells.Item(27, 2) "Cosmos"
Set-CellText $ws.Cells.Item(27, 3) "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-CellText $ws.Cells.Item(27, 4) "9.31"
Set-CellText $ws.Cells.Item(27, 5) "  +0.56%  "
Set-CellText $ws.Cells.Item(28, 5) "  -1.84%  "
Set-CellText $ws.Cells.Item(29, 4) "18.84"
Set-CellText $ws.Cells.Item(29, 5) "  -0.96%  "
Set-CellText $ws.Cells.Item(30, 4) "1.28"
Set-CellText $ws.Cells.Item(30, 5) "  -3.75%  "
Set-CellText $ws.Cells.Item(31, 4) "0.120"
Set-CellText $ws.Cells.Item(31, 5) "  +0.96%  "
Set-CellText $ws.Cells.Item(32, 4) "2.19"
Set-CellText $ws.Cells.Item(32, 5) "  +6.09%  "
Set-CellText $ws.Cells.Item(33, 5) "  -3.26%  "
Set-CellText $ws.Cells.Item(34, 5) "  +0.04%  "
Set-CellText $ws.Cells.Item(35, 5) "  -1.93%  "
Set-CellText $ws.Cells.Item(36, 5) "  +2.58%  "
Set-CellText $ws.Cells.Item(37, 4) "2.29"
Set-CellText $ws.Cells.Item(37, 5) "  -2.03%  "
Set-CellText $ws.Cells.Item(38, 4) "3.34"
Set-CellText $ws.Cells.Item(38, 5) "  +1.93%  "
Set-CellText $ws.Cells.Item(39, 5) "  +0.11%  "
Set-CellText $ws.Cells.Item(40, 4) "17.81"
Set-CellText $ws.Cells.Item(40, 5) "  +4.35%  "
Set-CellText $ws.Cells.Item(41, 4) "1.536.29"
Set-CellText $ws.Cells.Item(41, 5) "  +0.20%  "
Set-CellText $ws.Cells.Item(42, 5) "  +0.11%  "
Set-CellText $ws.Cells.Item(43, 4) "96.14"
Set-CellText $ws.Cells.Item(43, 5) "  -2.12%  "
Set-CellText $ws.Cells.Item(44, 5) "  -2.74%  "
Set-CellText $ws.Cells.Item(45, 4) "0.0915"
Set-CellText $ws.Cells.Item(45, 5) "  -0.95%  "
Set-CellText $ws.Cells.Item(46, 4) "1.11"
Set-CellText $ws.Cells.Item(46, 5) "  -2.24%  "
Set-CellText $ws.Cells.Item(47, 5) "  -2.74%  "
Set-CellText $ws.Cells.Item(48, 5) "  -1.32%  "
Set-CellText $ws.Cells.Item(49, 4) "7.17"
Set-CellText $ws.Cells.Item(49, 5) "  +0.72%  "
Set-CellText $ws.Cells.Item(50, 5) "  -0.47%  "
Set-CellText $ws.Cells.Item(51, 4) "2.226.80"
Set-CellText $ws.Cells.Item(51, 5) "  -1.04%  "
